$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messages")

# Day 1's image/audio links were re-uploaded; Excel's query refresh appended
# the "=FILE_ID" suffix that Google Drive adds to these two urls.
$c3 = $ws.Range("C3")
$d3 = $ws.Range("D3")
$d3.Value = $d3.Value() + "=FILE_ID"
$c3.Value = $c3.Value() + "=FILE_ID"

# Leave the active selection on C3 (image column) instead of D3.
$ws.Range("C3").Select() | Out-Null
